$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- PROTOTYPE section additions (Austin + Chris tasks for the prototype) ---

# Register the "Basic Networking, Basic Replication" string first so it lands
# in the shared-string table ahead of "PROTOTYPE" (matches authoring order).
$ws.Range("D14").Value = "Basic Networking, Basic Replication"

# Row 13: date entry (Aug 6, 2016), centered horizontally + vertically,
# using the same "d-mmm" date format as the existing date cell (B5).
$ws.Range("B13").NumberFormat = "d-mmm"
$ws.Range("B13").VerticalAlignment = -4108
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("B13").Value = 42588

# Row 12: "PROTOTYPE" section header, horizontally centered.
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("B12").Value = "PROTOTYPE"

# Row 14: Austin's prototype task.
$ws.Range("B14").Value = "Austin"

# Row 15: Chris's prototype task.
$ws.Range("B15").Value = "Chris"
$ws.Range("D15").Value = "Sounds"

# Move/restore the active selection to just below the new content.
$ws.Range("B16").Select()
